$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle header bands to the colored/filled header style (matching C4's style) ---
$ws.Range("C4").Copy()
$ws.Range("M4:P4").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C19:F19").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("H19:K19").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("M19:P19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in results for the 4/9-processor Tape Algorithm table (rows 6-10, columns N:P) ---
$ws.Range("N6").Value = "'0.0059"
$ws.Range("O6").Value = "'0.0017"
$ws.Range("P6").Value = "'3.3941"
$ws.Range("N7").Value = "'0.0673"
$ws.Range("O7").Value = "'0.0133"
$ws.Range("P7").Value = "'5.0667"
$ws.Range("N8").Value = "'0.5770"
$ws.Range("O8").Value = "'0.1260"
$ws.Range("P8").Value = "'4.5812"
$ws.Range("N9").Value = "'11.3943"
$ws.Range("O9").Value = "'0.9517"
$ws.Range("P9").Value = "'11.9721"
$ws.Range("N10").Value = "'150.5690"
$ws.Range("O10").Value = "'7.7621"
$ws.Range("P10").Value = "'19.3980"
$ws.Range("D21").Value = "'0.0025"
$ws.Range("E21").Value = "'0.0015"
$ws.Range("F21").Value = "'1.6885"
$ws.Range("I21").Value = "'0.0026"
$ws.Range("J21").Value = "'0.0017"
$ws.Range("K21").Value = "'1.5307"
$ws.Range("N21").Value = "'0.0026"
$ws.Range("O21").Value = "'0.0013"
$ws.Range("P21").Value = "'2.0102"
$ws.Range("D22").Value = "'0.0791"
$ws.Range("E22").Value = "'0.0259"
$ws.Range("F22").Value = "'3.0561"
$ws.Range("I22").Value = "'0.0796"
$ws.Range("J22").Value = "'0.0345"
$ws.Range("K22").Value = "'2.3089"
$ws.Range("N22").Value = "'0.0969"
$ws.Range("O22").Value = "'0.0274"
$ws.Range("P22").Value = "'3.5356"
$ws.Range("D23").Value = "'0.6746"
$ws.Range("E23").Value = "'0.1764"
$ws.Range("F23").Value = "'3.8235"
$ws.Range("I23").Value = "'0.6981"
$ws.Range("J23").Value = "'0.2170"
$ws.Range("K23").Value = "'3.2179"
$ws.Range("N23").Value = "'0.8964"
$ws.Range("O23").Value = "'0.1708"
$ws.Range("P23").Value = "'5.2485"
$ws.Range("D24").Value = "'2.7733"
$ws.Range("E24").Value = "'0.7017"
$ws.Range("F24").Value = "'3.9523"
$ws.Range("I24").Value = "'2.9752"
$ws.Range("J24").Value = "'0.8815"
$ws.Range("K24").Value = "'3.3751"
$ws.Range("N24").Value = "'3.4900"
$ws.Range("O24").Value = "'0.5009"
$ws.Range("P24").Value = "'6.9678"
$ws.Range("D25").Value = "'18.1305"
$ws.Range("E25").Value = "'5.6313"
$ws.Range("F25").Value = "'3.2196"
$ws.Range("I25").Value = "'17.8859"
$ws.Range("J25").Value = "'5.4933"
$ws.Range("K25").Value = "'3.2559"
$ws.Range("N25").Value = "'25.2653"
$ws.Range("O25").Value = "'2.2901"
$ws.Range("P25").Value = "'11.0323"
$ws.Range("C21").Value = 99
$ws.Range("H21").Value = 99
$ws.Range("M21").Value = 99
$ws.Range("C22").Value = 300
$ws.Range("H22").Value = 300
$ws.Range("M22").Value = 300
$ws.Range("C23").Value = 600
$ws.Range("H23").Value = 600
$ws.Range("M23").Value = 600
$ws.Range("C24").Value = 900
$ws.Range("H24").Value = 900
$ws.Range("M24").Value = 900
$ws.Range("C25").Value = 1500
$ws.Range("H25").Value = 1500
$ws.Range("M25").Value = 1500

# --- Update the selected cell to match the saved view state ---
$ws.Range("N28").Select()
